$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("U2").Value = 4.1
$ws.Range("V2").Value = 1.24
$ws.Range("N3").Value = 3.75
$ws.Range("Q3").Value = 3.5
$ws.Range("R3").Value = 1.31
$ws.Range("U3").Value = 9.4
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("O4").Value = 1.53
$ws.Range("P4").Value = 2.38
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.78
$ws.Range("U4").Value = 4.6
$ws.Range("V4").Value = 1.2
$ws.Range("W4").Value = 5.5
$ws.Range("X4").Value = 1.14
$ws.Range("AK4").Value = 21
$ws.Range("AO4").Value = 9.5
$ws.Range("G5").Value = 1.8
$ws.Range("J5").Value = 2.5
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("Q5").Value = 1.95
$ws.Range("R5").Value = 1.9
$ws.Range("AF5").Value = 13
$ws.Range("AQ5").Value = 51
$ws.Range("G6").Value = 2.63
$ws.Range("I6").Value = 3.2
$ws.Range("J6").Value = 3.6
$ws.Range("L6").Value = 4
$ws.Range("S6").Value = 3.4
$ws.Range("T6").Value = 1.33
$ws.Range("W6").Value = 7
$ws.Range("X6").Value = 1.1
$ws.Range("Y6").Value = 1.73
$ws.Range("Z6").Value = 2.08
$ws.Range("AD6").Value = 11
$ws.Range("AE6").Value = 12
$ws.Range("AH6").Value = 51
$ws.Range("AN6").Value = 13
$ws.Range("AP6").Value = 34
$ws.Range("G7").Value = 1.53
$ws.Range("H7").Value = 3.75
$ws.Range("I7").Value = 7
$ws.Range("J7").Value = 2.2
$ws.Range("K7").Value = 2.1
$ws.Range("L7").Value = 7
$ws.Range("N7").Value = 7.5
$ws.Range("O7").Value = 1.4
$ws.Range("P7").Value = 2.75
$ws.Range("S7").Value = 2.3
$ws.Range("T7").Value = 1.6
$ws.Range("AA7").Value = 2.38
$ws.Range("AB7").Value = 1.53
$ws.Range("AD7").Value = 6
$ws.Range("AF7").Value = 10
$ws.Range("AG7").Value = 15
$ws.Range("AI7").Value = 7.5
$ws.Range("AJ7").Value = 7.5
$ws.Range("AK7").Value = 23
$ws.Range("AM7").Value = 13
$ws.Range("AN7").Value = 34
$ws.Range("AP7").Value = 81
$ws.Range("AR7").Value = 67
$ws.Range("G8").Value = 1.9
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
$ws.Range("AA8").Value = 2.1
$ws.Range("AB8").Value = 1.67
$ws.Range("AI8").Value = 7
$ws.Range("AK8").Value = 19
$ws.Range("M9").Value = 1.11
$ws.Range("N9").Value = 6.5
$ws.Range("Q9").Value = 1.98
$ws.Range("R9").Value = 1.88
$ws.Range("S9").Value = 2.6
$ws.Range("T9").Value = 1.48
$ws.Range("U9").Value = 4
$ws.Range("V9").Value = 1.23
$ws.Range("G15").Value = 4.3
$ws.Range("H15").Value = 3.75
$ws.Range("I15").Value = 1.7
$ws.Range("J15").Value = 4.4
$ws.Range("K15").Value = 2.37
$ws.Range("L15").Value = 2.12
$ws.Range("P15").Value = 3.6
$ws.Range("AB15").Value = 2.05
$ws.Range("AC15").Value = 15
$ws.Range("AD15").Value = 28
$ws.Range("AF15").Value = 75
$ws.Range("AI15").Value = 13
$ws.Range("AJ15").Value = 7.6
$ws.Range("AK15").Value = 13.5
$ws.Range("AM15").Value = 8.25
$ws.Range("AO15").Value = 8
$ws.Range("AP15").Value = 13.5
$ws.Range("AR15").Value = 21
